# Generate Report for Handoff
# Adds a new file "9e17c9e5-2201-45a6-a90f-d80b36a112e3.md" as a new row to
# the Overview / zh-cn / de-de tables.

$wb = $excel.ActiveWorkbook

$hlColor = 15570276  # BGR for RGB(0x64,0x95,0xED) -> matches existing HyperLink font color FF6495ED
$dateFmt = "yyyy-mm-dd HH:mm:ss"

function Style-Hyperlink($cell) {
    $cell.Font.Underline = 2
    $cell.Font.Color = $hlColor
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 11
}

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$lo1 = $ws1.ListObjects.Item(1)
$lo1.ListRows.Add() | Out-Null

$ws1.Range("A8").Value = "9e17c9e5-2201-45a6-a90f-d80b36a112e3.md"
$ws1.Range("B8").Value = "e2e\9e17c9e5-2201-45a6-a90f-d80b36a112e3.md"
$ws1.Range("C8").Value = ".md"
$ws1.Range("E8").Value = "Ready for handoff"
$ws1.Range("F8").Value = "Ready for handoff"
$ws1.Range("G8").Value = "2016-11-14 17:46:43"
$ws1.Range("G8").NumberFormat = $dateFmt

$ws1.Hyperlinks.Add($ws1.Range("B8"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/90e23d7f8623667a77f103362734446b22d4ff6e/e2e/9e17c9e5-2201-45a6-a90f-d80b36a112e3.md", "", "", "e2e\9e17c9e5-2201-45a6-a90f-d80b36a112e3.md") | Out-Null
Style-Hyperlink $ws1.Range("B8")

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$lo2 = $ws2.ListObjects.Item(1)
$lo2.ListRows.Add() | Out-Null

$ws2.Range("A8").Value = "9e17c9e5-2201-45a6-a90f-d80b36a112e3.md"
$ws2.Range("B8").Value = ".md"
$ws2.Range("C8").Value = "Ready for handoff"
$ws2.Range("D8").Value = "e2e"
$ws2.Range("E8").Value = "ht"
$ws2.Range("F8").Value = "False"
$ws2.Range("G8").Value = "9e17c9e5-2201-45a6-a90f-d80b36a112e3.3152b83a2f40247f2188c99ab2de2c0f87dbbd86.zh-cn.xlf"
$ws2.Range("H8").Value = "2016-11-14 17:46:27"
$ws2.Range("H8").NumberFormat = $dateFmt
$ws2.Range("K8").Value = "0001-01-01 00:00:00"
$ws2.Range("K8").NumberFormat = $dateFmt
$ws2.Range("M8").Value = "True"
$ws2.Range("O8").Value = "False"

$ws2.Hyperlinks.Add($ws2.Range("A8"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/90e23d7f8623667a77f103362734446b22d4ff6e/e2e/9e17c9e5-2201-45a6-a90f-d80b36a112e3.md", "", "", "9e17c9e5-2201-45a6-a90f-d80b36a112e3.md") | Out-Null
Style-Hyperlink $ws2.Range("A8")

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$lo3 = $ws3.ListObjects.Item(1)
$lo3.ListRows.Add() | Out-Null

$ws3.Range("A8").Value = "9e17c9e5-2201-45a6-a90f-d80b36a112e3.md"
$ws3.Range("B8").Value = ".md"
$ws3.Range("C8").Value = "Ready for handoff"
$ws3.Range("D8").Value = "e2e"
$ws3.Range("E8").Value = "ht"
$ws3.Range("F8").Value = "False"
$ws3.Range("G8").Value = "9e17c9e5-2201-45a6-a90f-d80b36a112e3.3152b83a2f40247f2188c99ab2de2c0f87dbbd86.de-de.xlf"
$ws3.Range("H8").Value = "2016-11-14 17:46:43"
$ws3.Range("H8").NumberFormat = $dateFmt
$ws3.Range("K8").Value = "0001-01-01 00:00:00"
$ws3.Range("K8").NumberFormat = $dateFmt
$ws3.Range("M8").Value = "True"
$ws3.Range("O8").Value = "False"

$ws3.Hyperlinks.Add($ws3.Range("A8"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/90e23d7f8623667a77f103362734446b22d4ff6e/e2e/9e17c9e5-2201-45a6-a90f-d80b36a112e3.md", "", "", "9e17c9e5-2201-45a6-a90f-d80b36a112e3.md") | Out-Null
Style-Hyperlink $ws3.Range("A8")

Write-Host "Added handoff row for 9e17c9e5-2201-45a6-a90f-d80b36a112e3.md"
